# Trade #4 closed at 2026-02-17 19:55:59 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1399.78
$summary.Range("B4").Value = -0.22
$summary.Range("B5").Value = -1.1
$summary.Range("B6").Value = 4
$summary.Range("B7").Value = 1
$summary.Range("B9").Value = 25

# --- Strategy Status sheet (MarketMaking row = row 5) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 99.78
$status.Range("D5").Value = 4
$status.Range("E5").Value = -0.22
$status.Range("F5").Value = -0.22
$status.Range("G5").Value = 25

# --- New trade row to append to "All Trades" and "MarketMaking" sheets ---
$tradeRow = @{
    A = 4
    B = "2026-02-17"
    C = "19:55:52"
    D = "MarketMaking"
    E = "UP"
    F = 0.59
    G = 0.6
    H = "CLOSED"
    I = 1.6949
    J = 0.01
    K = 99.78
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.14
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A5").Value = $tradeRow.A
    $ws.Range("B5").NumberFormat = "@"
    $ws.Range("B5").Value = $tradeRow.B
    $ws.Range("C5").Value = $tradeRow.C
    $ws.Range("D5").Value = $tradeRow.D
    $ws.Range("E5").Value = $tradeRow.E
    $ws.Range("F5").Value = $tradeRow.F
    $ws.Range("G5").Value = $tradeRow.G
    $ws.Range("H5").Value = $tradeRow.H
    $ws.Range("I5").Value = $tradeRow.I
    $ws.Range("J5").Value = $tradeRow.J
    $ws.Range("K5").Value = $tradeRow.K
    $ws.Range("L5").Value = $tradeRow.L
    $ws.Range("M5").Value = $tradeRow.M
    $ws.Range("N5").Value = $tradeRow.N
    $ws.Range("O5").Value = $tradeRow.O
    $ws.Range("P5").Value = $tradeRow.P
    $ws.Range("Q5").Value = $tradeRow.Q
}
